$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 656.4666999999999
$ws.Range("I80").Value = 528.44446
$ws.Range("J80").Value = 848.5
$ws.Range("K80").Value = 1585.33338
$ws.Range("L80").Value = 2545.5
$ws.Range("M80").Value = -587.33338
$ws.Range("N80").Value = -4541.5
# Row 83
$ws.Range("H83").Value = 656.4666999999999
$ws.Range("I83").Value = 528.44446
$ws.Range("J83").Value = 848.5
$ws.Range("K83").Value = 4756.00014
$ws.Range("L83").Value = 7636.5
$ws.Range("M83").Value = 235.9998599999999
$ws.Range("N83").Value = -17620.5
# Row 86
$ws.Range("H86").Value = 3376.7778
$ws.Range("I86").Value = 3315.1667
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 3315.1667
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2192.1667
# Row 89
$ws.Range("H89").Value = 3376.7778
$ws.Range("I89").Value = 3315.1667
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 16575.8335
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -10959.8335
# Row 111
$ws.Range("H111").Value = 4000
$ws.Range("I111").Value = 4000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 12000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -8933
$ws.Range("N111").ClearContents()
# Row 113
$ws.Range("H113").Value = 5660.5264
$ws.Range("I113").Value = 4666.8887
$ws.Range("J113").Value = 6554.8
$ws.Range("K113").Value = 4666.8887
$ws.Range("L113").Value = 6554.8
$ws.Range("M113").Value = -1412.8887
$ws.Range("N113").Value = -13062.8
# Row 116
$ws.Range("H116").Value = 56415
$ws.Range("I116").Value = 40000
$ws.Range("J116").Value = 58760
$ws.Range("K116").Value = 40000
$ws.Range("L116").Value = 58760
$ws.Range("M116").Value = -36558
# Row 132
$ws.Range("H132").Value = 4017.1072
$ws.Range("I132").Value = 1924.9546
$ws.Range("J132").Value = 11688.333
$ws.Range("K132").Value = 5774.8638
$ws.Range("L132").Value = 35064.999
$ws.Range("M132").Value = -3244.8638
# Row 137
$ws.Range("H137").Value = 3587.5217
$ws.Range("I137").Value = 2979.5
$ws.Range("J137").Value = 3911.8
$ws.Range("K137").Value = 8938.5
$ws.Range("L137").Value = 11735.4
$ws.Range("M137").Value = -6388.5
$ws.Range("N137").Value = -16835.4

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3301070.8
$ws.Range("I32").Value = 570281.25
$ws.Range("J32").Value = 17930302
$ws.Range("K32").Value = 570281.25
$ws.Range("L32").Value = 17930302
$ws.Range("M32").Value = -569994.25
# Row 74
$ws.Range("H74").Value = 2811.3076
$ws.Range("I74").Value = 1920.2
$ws.Range("J74").Value = 5781.6665
$ws.Range("K74").Value = 1920.2
$ws.Range("L74").Value = 5781.6665
$ws.Range("M74").Value = -1046.2
# Row 77
$ws.Range("H77").Value = 2811.3076
$ws.Range("I77").Value = 1920.2
$ws.Range("J77").Value = 5781.6665
$ws.Range("K77").Value = 9601
$ws.Range("L77").Value = 28908.3325
$ws.Range("M77").Value = -5233
# Row 102
$ws.Range("H102").Value = 2016.8
$ws.Range("I102").Value = 2016.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2016.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -394.8
$ws.Range("N102").ClearContents()
# Row 122
$ws.Range("H122").Value = 1740
$ws.Range("I122").Value = 1551.4286
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4654.2858
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2204.2858
$ws.Range("N122").Value = -12100
# Row 132
$ws.Range("H132").Value = 2694.3784
$ws.Range("I132").Value = 2694.3784
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8083.135200000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5553.135200000001

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2392.48
$ws.Range("I31").Value = 1535.6875
$ws.Range("J31").Value = 3915.6667
$ws.Range("K31").Value = 1535.6875
$ws.Range("L31").Value = 3915.6667
$ws.Range("M31").Value = -1240.6875
# Row 34
$ws.Range("H34").Value = 2392.48
$ws.Range("I34").Value = 1535.6875
$ws.Range("J34").Value = 3915.6667
$ws.Range("K34").Value = 1535.6875
$ws.Range("L34").Value = 3915.6667
$ws.Range("M34").Value = -1333.6875

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 601.1818
$ws.Range("I12").Value = 686
$ws.Range("J12").Value = 552.7143
$ws.Range("K12").Value = 2058
$ws.Range("L12").Value = 1658.1429
$ws.Range("M12").Value = -1885
$ws.Range("N12").Value = -2004.1429
# Row 14
$ws.Range("H14").Value = 279.9091
$ws.Range("I14").Value = 279.9091
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 839.7273
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -666.7273
# Row 41
$ws.Range("H41").Value = 3595
$ws.Range("I41").Value = 286.33334
$ws.Range("J41").Value = 5249.3335
$ws.Range("K41").Value = 859.0000200000001
$ws.Range("L41").Value = 15748.0005
$ws.Range("M41").Value = -521.0000200000001
$ws.Range("N41").Value = -16424.0005
# Row 98
$ws.Range("H98").Value = 208.33333
$ws.Range("I98").Value = 195
$ws.Range("J98").Value = 225
$ws.Range("K98").Value = 585
$ws.Range("L98").Value = 675
$ws.Range("M98").Value = 913
$ws.Range("N98").Value = -3671
# Row 122
$ws.Range("H122").Value = 447.45456
$ws.Range("I122").Value = 304
$ws.Range("J122").Value = 461.8
$ws.Range("K122").Value = 2736
$ws.Range("L122").Value = 4156.2
$ws.Range("M122").Value = -286
$ws.Range("N122").Value = -9056.200000000001
# Row 138
$ws.Range("H138").Value = 540
$ws.Range("I138").Value = 540
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 1620
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 3520
# Row 139
$ws.Range("H139").Value = 206350.6
$ws.Range("I139").Value = 999999
$ws.Range("J139").Value = 7938.5
$ws.Range("K139").Value = 2999997
$ws.Range("L139").Value = 23815.5
$ws.Range("M139").Value = -2994857

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 72151180
$ws.Range("I70").Value = 5258.5713
$ws.Range("J70").Value = 173155490
$ws.Range("K70").Value = 5258.5713
$ws.Range("L70").Value = 173155490
$ws.Range("M70").Value = -4988.5713
# Row 73
$ws.Range("H73").Value = 72151180
$ws.Range("I73").Value = 5258.5713
$ws.Range("J73").Value = 173155490
$ws.Range("K73").Value = 5258.5713
$ws.Range("L73").Value = 173155490
$ws.Range("M73").Value = -4322.5713
# Row 80
$ws.Range("H80").Value = 9674.143
$ws.Range("I80").Value = 12838.667
$ws.Range("J80").Value = 3978
$ws.Range("K80").Value = 12838.667
$ws.Range("L80").Value = 3978
$ws.Range("M80").Value = -11840.667
$ws.Range("N80").Value = -5974
# Row 83
$ws.Range("H83").Value = 9674.143
$ws.Range("I83").Value = 12838.667
$ws.Range("J83").Value = 3978
$ws.Range("K83").Value = 64193.335
$ws.Range("L83").Value = 19890
$ws.Range("M83").Value = -59201.335
$ws.Range("N83").Value = -29874
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
# Row 141
$ws.Range("H141").Value = 93999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 93999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 93999.5
$ws.Range("N141").Value = -104359.5

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6943.1
$ws.Range("I40").Value = 6492.7144
$ws.Range("J40").Value = 7994
$ws.Range("K40").Value = 6492.7144
$ws.Range("L40").Value = 7994
$ws.Range("M40").Value = -6356.7144
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
# Row 51
$ws.Range("H51").Value = 19116.666
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19116.666
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19116.666
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -20072.666
# Row 55
$ws.Range("H55").Value = 407.72223
$ws.Range("I55").Value = 359.8
$ws.Range("J55").Value = 467.625
$ws.Range("K55").Value = 359.8
$ws.Range("L55").Value = 467.625
$ws.Range("M55").Value = -186.8
# Row 82
$ws.Range("H82").Value = 5786
$ws.Range("I82").Value = 6580.4
$ws.Range("J82").Value = 3800
$ws.Range("K82").Value = 6580.4
$ws.Range("L82").Value = 3800
$ws.Range("M82").Value = -6219.4
$ws.Range("N82").Value = -4522
# Row 85
$ws.Range("H85").Value = 5786
$ws.Range("I85").Value = 6580.4
$ws.Range("J85").Value = 3800
$ws.Range("K85").Value = 6580.4
$ws.Range("L85").Value = 3800
$ws.Range("M85").Value = -5332.4
$ws.Range("N85").Value = -6296
# Row 122
$ws.Range("H122").Value = 4816.3335
$ws.Range("I122").Value = 4299.6665
$ws.Range("J122").Value = 5333
$ws.Range("K122").Value = 12898.9995
$ws.Range("L122").Value = 15999
$ws.Range("M122").Value = -10448.9995
# Row 132
$ws.Range("H132").Value = 4048.9
$ws.Range("I132").Value = 2915.6667
$ws.Range("J132").Value = 5748.75
$ws.Range("K132").Value = 8747.000100000001
$ws.Range("L132").Value = 17246.25
$ws.Range("M132").Value = -6217.000100000001
$ws.Range("N132").Value = -22306.25
# Row 136
$ws.Range("H136").Value = 5248.25
$ws.Range("I136").Value = 4998.2
$ws.Range("J136").Value = 5665
$ws.Range("K136").Value = 14994.6
$ws.Range("L136").Value = 16995
$ws.Range("M136").Value = -12444.6
$ws.Range("N136").Value = -22095

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 25666.334
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 25666.334
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 51332.668
$ws.Range("N81").Value = -53454.668
# Row 84
$ws.Range("H84").Value = 25666.334
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 25666.334
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 256663.34
$ws.Range("N84").Value = -267271.34
# Row 107
$ws.Range("H107").Value = 45507090
$ws.Range("I107").Value = 1283.8334
$ws.Range("J107").Value = 100114050
$ws.Range("K107").Value = 3851.5002
$ws.Range("L107").Value = 300342150
$ws.Range("M107").Value = -1931.5002
$ws.Range("N107").Value = -300345990
# Row 113
$ws.Range("H113").Value = 1568.7646
$ws.Range("I113").Value = 1074
$ws.Range("J113").Value = 1915.1
$ws.Range("K113").Value = 3222
$ws.Range("L113").Value = 5745.299999999999
$ws.Range("M113").Value = -1052
$ws.Range("N113").Value = -10085.3
# Row 136
$ws.Range("H136").Value = 5873.9414
$ws.Range("I136").Value = 4178.3335
$ws.Range("J136").Value = 12414.143
$ws.Range("K136").Value = 12535.0005
$ws.Range("L136").Value = 37242.429
$ws.Range("M136").Value = -9985.000499999998
# Row 140
$ws.Range("H140").Value = 79391.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 79391.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 79391.8
$ws.Range("N140").Value = -89751.8
